$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.204.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.800.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5269'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3818'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07997'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.321'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.63'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.804.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.312'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001096'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06615'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.973'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.245.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.239'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.007.99'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.388'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1095'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.062'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.66%  '
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.533'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07290'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.77%  '
$ws.Range('E36').Value = '  +9.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2165'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.832'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02310'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.060'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6195'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.369'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5980'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.758'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.202'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.928'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06833'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.86'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.66%  '
